# Apply the "I0"/"IF" column addition described in the diff:
# two new header cells (I1, J1) styled like the existing headers,
# plus nine rows of numeric data in columns I and J.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style from existing header cell (H1) to new header cells first,
# then set their values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 11
$ws.Range("J4").Value = 11

$ws.Range("I5").Value = 8
$ws.Range("J5").Value = 8

$ws.Range("I6").Value = 9
$ws.Range("J6").Value = 9

$ws.Range("I7").Value = 9
$ws.Range("J7").Value = 9

$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 1

$ws.Range("I9").Value = 9
$ws.Range("J9").Value = 9
